# OPtoCUMappings.xlsx - add "Hex" column (L) computing DEC2HEX of the
# existing "Value" column (K), which itself becomes a computed bitmask
# formula instead of a hand-typed literal. Also fixes two stray bits in
# the MULT row (C5/D5) that were left over from copy/paste.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New "Hex" header in L2
# ---------------------------------------------------------------------
$ws.Range("L2").Value = "Hex"

# ---------------------------------------------------------------------
# 2. Fix the MULT row (row 5): RegDst/RegWrite bits were wrongly 1/1,
#    should be 0/0 like the other R-type math ops that don't write back
#    a register the same way.
# ---------------------------------------------------------------------
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

# ---------------------------------------------------------------------
# 3. Column K (rows 3-20): replace literal values with the bitmask
#    formula that derives "Value" from the individual control-signal
#    columns C:J.
# ---------------------------------------------------------------------
for ($row = 3; $row -le 20; $row++) {
    $ws.Cells.Item($row, 11).Formula = "=J$row*128+I$row*64+H$row*32+G$row*16+F$row*8+E$row*4+D$row*2+C$row"
}

# Rows 19 and 20 don't have any control-signal bits filled in, so the
# formula there would show a stray 0 - clear them back out but keep the
# (now fill-formatted) cell alive, matching the rest of column K.
$ws.Range("K19").ClearContents()
$ws.Range("K20").ClearContents()
$ws.Range("K18").Copy()
$ws.Range("K19:K20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Column L (rows 3-18): new "Hex" column, converts the decimal Value
#    in K to its hex-string representation.
# ---------------------------------------------------------------------
for ($row = 3; $row -le 18; $row++) {
    $ws.Cells.Item($row, 12).Formula = "=DEC2HEX(K$row)"
}

# The old L13:L16 cells had a leftover yellow highlight fill with no
# content - now that L has real data in every row, drop that formatting
# so the column looks consistent.
$ws.Range("L13:L16").ClearFormats()

# ---------------------------------------------------------------------
# 5. Selection cosmetic change, matches the author's saved cursor spot.
# ---------------------------------------------------------------------
$ws.Range("N3").Select()
